# Pipeline update: move the "DGJ-DAVI" participant row (row 54) down to
# just above the final row ("luck"), shifting the rows in between up by
# one. Net effect: row 54's original data is removed, rows 55-70 shift up
# to become rows 54-69, and a new row 70 is created holding DGJ-DAVI's
# original stats; row 71 ("luck") is left untouched at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the DGJ-DAVI row; everything below shifts up one row.
$ws.Rows(54).Delete()

# Make room for DGJ-DAVI again immediately above the last row ("luck",
# now sitting at row 70 after the shift).
$ws.Rows(70).Insert()

# Re-insert DGJ-DAVI's original data at the new row 70.
$ws.Range("A70").Value = "DGJ-DAVI"
$ws.Range("B70").Value = "Ok"
$ws.Range("C70").Value = "Guerra Atual"
$ws.Range("D70").Value = 16
$ws.Range("E70").Value = 16
$ws.Range("F70").Value = 16
$ws.Range("G70").Value = 12
$ws.Range("H70").Value = 15
